$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan")

# Column width changes: B 55.71 -> 51.71, E 45.71 -> 55.71
# (values chosen land on the nearest width this engine's pixel-quantized
# ColumnWidth setter can represent, closest to the target char-widths)
$ws.Columns.Item(2).ColumnWidth = 50.8
$ws.Columns.Item(5).ColumnWidth = 54.8

# New schedule grid for rows 2-13, columns B-F
$grid = @{
    2  = @{ B = "{0: sala nr 4 | Lena Kowalska | Język angielski}";      C = "{0: sala nr 11 | Katarzyna Mazur | Fizyka}";            D = "{0: sala nr 11 | Dominik Kaczor | Informatyka}";       E = "{}";                                                    F = "{}" }
    3  = @{ B = "{0: sala nr 4 | Mateusz Kowalski | Język niemiecki}";   C = "{0: sala nr 2 | Mateusz Kowalski | Język niemiecki}";  D = "{0: sala nr 10 | Paweł Lewandowski | Matematyka}";     E = "{}";                                                    F = "{}" }
    4  = @{ B = "{0: sala nr 3 | Karolina Kamińska | Chemia}";           C = "{0: sala nr 1 | Jan Nowak | Język polski}";            D = "{0: sala nr 2 | Paweł Lewandowski | Matematyka}";      E = "{}";                                                    F = "{}" }
    5  = @{ B = "{0: sala nr 8 | Dominik Kaczor | Informatyka}";         C = "{0: sala nr 3 | Zofia Wiśniewska | Wychowanie fizyczne}"; D = "{0: sala nr 11 | Paweł Lewandowski | Matematyka}";  E = "{}";                                                    F = "{}" }
    6  = @{ B = "{0: sala nr 7 | Katarzyna Mazur | Fizyka}";             C = "{0: sala nr 5 | Natalia Szymańska | Geografia}";       D = "{0: sala nr 9 | Piotr Wójcik | Biologia}";             E = "{0: sala nr 9 | Paweł Lewandowski | Matematyka}";      F = "{}" }
    7  = @{ B = "{0: sala nr 7 | Piotr Wójcik | Biologia}";              C = "{0: sala nr 7 | Dominik Kaczor | Informatyka}";        D = "{}";                                                    E = "{0: sala nr 2 | Natalia Szymańska | Geografia}";       F = "{}" }
    8  = @{ B = "{}";                                                    C = "{0: sala nr 3 | Jan Nowak | Język polski}";            D = "{}";                                                    E = "{}";                                                    F = "{0: sala nr 5 | Paweł Lewandowski | Matematyka}" }
    9  = @{ B = "{}";                                                    C = "{}";                                                    D = "{}";                                                    E = "{0: sala nr 4 | Dominik Kaczor | Informatyka}";        F = "{0: sala nr 9 | Lena Kowalska | Język angielski}" }
    10 = @{ B = "{}";                                                    C = "{}";                                                    D = "{}";                                                    E = "{0: sala nr 2 | Zofia Wiśniewska | Wychowanie fizyczne}"; F = "{0: sala nr 1 | Zofia Wiśniewska | Wychowanie fizyczne}" }
    11 = @{ B = "{}";                                                    C = "{}";                                                    D = "{}";                                                    E = "{0: sala nr 5 | Lena Kowalska | Język angielski}";     F = "{}" }
    12 = @{ B = "{}";                                                    C = "{}";                                                    D = "{}";                                                    E = "{0: sala nr 11 | Jan Nowak | Język polski}";           F = "{0: sala nr 2 | Karolina Kamińska | Chemia}" }
    13 = @{ B = "{}";                                                    C = "{}";                                                    D = "{}";                                                    E = "{0: sala nr 1 | Katarzyna Mazur | Fizyka}";            F = "{0: sala nr 3 | Dominik Kaczor | Informatyka}" }
}

foreach ($row in $grid.Keys) {
    $cols = $grid[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
